# "Generate Report for Handoff"
#
# The handoff/handback filenames change from the
# fead2b92-59e1-4e44-a218-a1b33fb17d10 GUID to a fresh
# bb3882a2-41ea-46eb-9be9-889a5abac812 GUID, and the associated timestamps
# advance a few seconds/minutes to reflect the newly generated report.

$wb = $excel.ActiveWorkbook

$oldGuid = "fead2b92-59e1-4e44-a218-a1b33fb17d10"
$newGuid = "bb3882a2-41ea-46eb-9be9-889a5abac812"

$oldMdName = "$oldGuid.md"
$newMdName = "$newGuid.md"

$newZhXlf = "$newGuid.a6b64076f672c31e6479f08778d69ca41daae540.zh-cn.xlf"
$newDeXlf = "$newGuid.a6b64076f672c31e6479f08778d69ca41daae540.de-de.xlf"

# Hyperlinks on all three sheets point at the same (unchanged) commit URL for
# the old handoff markdown file -- only the visible display text (and the
# backing cell text) move to the new filename.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02389192015ecf18eae8ee25be28a035eb1e4dc6/e2e/$oldMdName"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName

$overviewDisplay = "e2e\$newMdName"
$wsOverview.Range("B2").Value = $overviewDisplay
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", $overviewDisplay) | Out-Null

$wsOverview.Range("G2").Value = "2016-12-16 09:09:24"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", $newMdName) | Out-Null

$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-12-16 09:09:11"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", $newMdName) | Out-Null

$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-12-16 09:09:24"
